$d = $word.ActiveDocument

# --- New Heading1 paragraph, appended right after the last existing
#     paragraph ("When the player gets to the end goal ...") ---
$lastPara = $d.Paragraphs.Last
$rng = $lastPara.Range
$rng.Collapse(0)                 # wdCollapseEnd -> end of the last paragraph
$rng.InsertParagraphAfter()
$rng.Collapse(0)

$headingPara = $d.Paragraphs.Last
$headingPara.Style = "Heading1"
$headingPara.Range.Text = "Onion Design of Project Boost"

# --- Four new numbered ListParagraph items under the new heading ---
$items = @("Player movement", "Obstacle collision", "Progressing through levels", "Fuel mechanic")

$prevRange = $headingPara.Range
$prevRange.Collapse(0)           # wdCollapseEnd -> end of the heading paragraph

$firstItemStart = -1
for ($i = 0; $i -lt $items.Count; $i++) {
    $prevRange.InsertParagraphAfter()
    $prevRange.Collapse(0)

    $itemPara = $d.Paragraphs.Last
    $itemPara.Style = "ListParagraph"
    $itemPara.Range.Text = $items[$i]

    if ($firstItemStart -eq -1) {
        $firstItemStart = $itemPara.Range.Start
    }

    $prevRange = $itemPara.Range
    $prevRange.Collapse(0)
}

# Apply numbered-list formatting to all four new items at once so they
# share a single new numbering definition (one new numId), the way Word
# does when you select several paragraphs and click "Numbering".
$lastItemEnd = $d.Paragraphs.Last.Range.End
$listRange = $d.Range($firstItemStart, $lastItemEnd)
$listRange.ListFormat.ApplyNumberDefault()

Write-Output ("Paragraphs: " + $d.Paragraphs.Count)
